# Corrected excel sheets for application fix issues
#
# 1. "Summary" sheet: selection moves from C5 to A7:XFD15.
# 2. "Repayment schedule" sheet: selection moves from F8 to A9:XFD9, and a
#    previously-missing "Outstanding Interest" column (O) is populated for
#    every data row (and a trailing blank P2 cell is restored), reusing the
#    existing formatting already present on the neighbouring "Late" (N)
#    column so no new cell styles are introduced.

$wb = $excel.ActiveWorkbook

# Remember whatever sheet is active right now so we can restore it once
# we're done poking at the other sheets (selecting a range on a sheet makes
# Excel activate that sheet first, same as a user clicking on it).
$originalSheet = $wb.ActiveSheet

# --- Summary sheet -----------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Select()
$summary.Range("A7:XFD15").Select()

# --- Repayment schedule sheet -------------------------------------------
$repay = $wb.Worksheets.Item("Repayment schedule")
$repay.Select()

# Restore the blank, formatted P2 cell (copy format from O2, which already
# carries the correct style).
$repay.Cells.Item(2, 15).Copy()
$repay.Cells.Item(2, 16).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in column O ("Outstanding" interest) for rows 3-8 with 0, copying
# the formatting from column N (the adjacent "Late" column) so the new
# cells pick up the same style already used throughout the sheet instead
# of registering a brand-new one.
for ($r = 3; $r -le 8; $r++) {
    $repay.Cells.Item($r, 14).Copy()
    $repay.Cells.Item($r, 15).PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $repay.Cells.Item($r, 15).Value = 0
}

$repay.Range("A9:XFD9").Select()

# Put the view back the way it was.
$originalSheet.Select()
